$wb = $excel.ActiveWorkbook

# Helper: write a value into a Range while forcing Excel to keep it as TEXT
# (shared string) even though it "looks like" a number. We briefly switch
# the cell to a text number format so Excel doesn't auto-coerce the value,
# then restore the "Normal" style so no stray formatting is left behind.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Sheet: Restricciones_del_follower (numeric expressions / values, stored as text) ---
$ws = $wb.Worksheets.Item(3)

Set-TextValue $ws.Range("A2") "4.366913451651778 - 0.7574285534004321y_1 + 0.4286380056576421y_2"
Set-TextValue $ws.Range("B2") "-4.366913451651778"
Set-TextValue $ws.Range("D2") "0.03"
Set-TextValue $ws.Range("E2") "3.3000000000000003"
Set-TextValue $ws.Range("F2") "0"

Set-TextValue $ws.Range("A3") "3.303086074127693 + 0.13373800663553037y_1 - 1.5307358653210148y_2"
Set-TextValue $ws.Range("B3") "-7.303086074127693"
Set-TextValue $ws.Range("D3") "0.85"
Set-TextValue $ws.Range("E3") "-0.8"
Set-TextValue $ws.Range("F3") "-0.5"

Set-TextValue $ws.Range("A4") "5.75309525256657 - 2x + 0.22725019028865578y_1 + 2.634505660289932y_2"
Set-TextValue $ws.Range("B4") "-21.75309525256657"
Set-TextValue $ws.Range("D4") "0.73"
Set-TextValue $ws.Range("E4") "0"
Set-TextValue $ws.Range("F4") "3.5"

Set-TextValue $ws.Range("A5") "-64.32691345165179 + 8x + 0.7574285534004321y_1 - 0.4286380056576421y_2"
Set-TextValue $ws.Range("B5") "15.566913451651782"
Set-TextValue $ws.Range("D5") "0.02"
Set-TextValue $ws.Range("E5") "1.1"
Set-TextValue $ws.Range("F5") "0"

Set-TextValue $ws.Range("A6") "6.602172717320021 - 2x - 1.0159046693889429y_1 + 1.7389542990459823y_2"
Set-TextValue $ws.Range("B6") "-5.397827282679979"
Set-TextValue $ws.Range("D6") "0.72"
Set-TextValue $ws.Range("E6") "3.1"
Set-TextValue $ws.Range("F6") "0.4"

# --- Sheet: Punto_modificado (x, y_1, y_2 stored as text) ---
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("A2") "7.4"
Set-TextValue $ws.Range("B2") "7.35"
Set-TextValue $ws.Range("C2") "2.8"

# --- Sheet: Vector_bf (stored as text) ---
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "1.4594557029431234"
Set-TextValue $ws.Range("A3") "-1.8783971218584716"

# --- Sheet: Vector_BF (stored as text) ---
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "-1.6000000000000005"
Set-TextValue $ws.Range("A3") "7.922637697895098"
Set-TextValue $ws.Range("A4") "-9.558350631746169"

# --- Sheet: Vector_Alpha (genuine numbers) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 0.19446458492740665
$ws.Range("A3").Value = 0.34363035313024864
